{"js": "// The document had two spots where Word had split one logical sentence\n// into two separate runs (likely from an earlier autocorrect / typing\n// edit). This change re-joins each pair of runs into a single run with\n// the same (bold) formatting, without altering the visible text.\n\nconst body = context.document.body;\n\n// --- Change 1: \"...roll numbe\" + \"r \" -> \"...roll number \" -----------\nconst target1 =\n  \"CREATE TABLE Student (Username varchar (20), email varchar (20), \" +\n  \"roll number \";\nconst results1 = body.search(target1, { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  results1.items[0].insertText(target1, Word.InsertLocation.replace);\n}\n\n// --- Change 2: \"...Insert into Student \" + \"values (...);  \" --------\nconst target2 =\n  \"Insert into Student values (\\u2018CD\\u2019, \\u2018CD@gmail.com\\u2019,\" +\n  \"\\u20198\\u2019,\\u201905\\u2019); Insert into Student values \" +\n  \"(\\u2018MN\\u2019, \\u2018MN@gmail.com\\u2019,\\u20190\\u2019,\\u201909\\u2019);  \";\nconst results2 = body.search(target2, { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  results2.items[0].insertText(target2, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document had two spots where one logical bold sentence had been\n# split across two separate <w:r> runs (same run formatting on both\n# halves). This re-joins each pair of runs into a single run with the\n# same (bold) text, leaving the visible text unchanged.\n#\n# Find.Execute positional args used below:\n#   (FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#    MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# wdFindContinue = 1 (Wrap), wdReplaceAll = 2 (Replace)\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"...roll numbe\" + \"r \" -> single run \"...roll number \" ---\n$target1 = \"CREATE TABLE Student (Username varchar (20), email varchar (20), roll number \"\n$range1 = $d.Content\n$range1.Find.ClearFormatting()\n$range1.Find.Replacement.ClearFormatting()\n$range1.Find.Execute($target1, $true, $false, $false, $false, $false, $true, 1, $false, $target1, 2)\n\n# --- Change 2: \"...Insert into Student \" + \"values (...);  \" -> single run ---\n$target2 = \"Insert into Student values (\" + [char]0x2018 + \"CD\" + [char]0x2019 + \", \" + [char]0x2018 + \"CD@gmail.com\" + [char]0x2019 + \",\" + [char]0x2019 + \"8\" + [char]0x2019 + \",\" + [char]0x2019 + \"05\" + [char]0x2019 + \"); Insert into Student values (\" + [char]0x2018 + \"MN\" + [char]0x2019 + \", \" + [char]0x2018 + \"MN@gmail.com\" + [char]0x2019 + \",\" + [char]0x2019 + \"0\" + [char]0x2019 + \",\" + [char]0x2019 + \"09\" + [char]0x2019 + \");  \"\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Replacement.ClearFormatting()\n$range2.Find.Execute($target2, $true, $false, $false, $false, $false, $true, 1, $false, $target2, 2)\n"}
